# Insert a new weekly data row for
# "Feria Lagunitas de Puerto Montt - Cilantro" right before the current
# row 72, shifting all subsequent rows (72-182) down by one (to 73-183).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 72..182 down to 73..183 and leave a blank row 72 to fill in.
$ws.Rows(72).Insert()

# Populate the newly inserted row 72 with the new week's data.
$ws.Cells.Item(72, 1).Value = 4
$ws.Cells.Item(72, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(72, 3).Value = "Los Lagos"
$ws.Cells.Item(72, 4).Value = 44495
$ws.Cells.Item(72, 5).Value = 10
$ws.Cells.Item(72, 6).Value = 100112040
$ws.Cells.Item(72, 7).Value = "Cilantro"
$ws.Cells.Item(72, 8).Value = "Sin especificar"
$ws.Cells.Item(72, 9).Value = "Primera"
$ws.Cells.Item(72, 10).Value = 300
$ws.Cells.Item(72, 11).Value = 10000
$ws.Cells.Item(72, 12).Value = 10000
$ws.Cells.Item(72, 13).Value = 10000
$ws.Cells.Item(72, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(72, 15).Value = "Región Metropolitana"
$ws.Cells.Item(72, 16).Value = 278
$ws.Cells.Item(72, 17).Value = 36
$ws.Cells.Item(72, 18).Value = "Hortaliza"

# Match the date-formatted style used by the rest of column D.
$ws.Cells.Item(72, 4).NumberFormat = $ws.Cells.Item(73, 4).NumberFormat
